$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1066.827089534028
$ws.Range("C2").Value = 1511.213283282402
$ws.Range("D2").Value = 1838.018729192084
$ws.Range("E2").Value = 2010.631070832981

$ws.Range("B3").Value = 1096.168533900459
$ws.Range("C3").Value = 1528.33010704203
$ws.Range("D3").Value = 1847.781862906447
$ws.Range("E3").Value = 2018.735400027925

$ws.Range("B4").Value = 1121.183306044844
$ws.Range("C4").Value = 1540.626423965915
$ws.Range("D4").Value = 1850.061862514113
$ws.Range("E4").Value = 2017.413604081494

$ws.Range("B5").Value = 1125.261196947115
$ws.Range("C5").Value = 1546.222606126279
$ws.Range("D5").Value = 1856.916359620758
$ws.Range("E5").Value = 2025.331113879538

$ws.Range("B6").Value = 1103.837453310326
$ws.Range("C6").Value = 1536.915018266464
$ws.Range("D6").Value = 1855.084687763432
$ws.Range("E6").Value = 2024.724179085728

$ws.Range("B7").Value = 1154.182764083716
$ws.Range("C7").Value = 1564.451171895376
$ws.Range("D7").Value = 1866.352806353311
$ws.Range("E7").Value = 2030.493497846352

$ws.Range("B8").Value = 1048.501688302613
$ws.Range("C8").Value = 1485.222641346011
$ws.Range("D8").Value = 1810.082880098869
$ws.Range("E8").Value = 1982.384839167365

$ws.Range("B9").Value = 1134.340219809489
$ws.Range("C9").Value = 1548.913673198975
$ws.Range("D9").Value = 1854.667135943324
$ws.Range("E9").Value = 2019.818827124698

$ws.Range("B10").Value = 1228.166133370289
$ws.Range("C10").Value = 1647.746076806086
$ws.Range("D10").Value = 1939.815918085547
$ws.Range("E10").Value = 2038.806492436251

$ws.Range("B11").Value = 1294.365777851517
$ws.Range("C11").Value = 1678.467432057428
$ws.Range("D11").Value = 1953.123059391344
$ws.Range("E11").Value = 2047.386004418435

$ws.Range("B12").Value = 1172.687020782348
$ws.Range("C12").Value = 1616.42836104628
$ws.Range("D12").Value = 1925.405581190644
$ws.Range("E12").Value = 2028.861029587857

$ws.Range("B13").Value = 1219.485333534142
$ws.Range("C13").Value = 1641.824615553334
$ws.Range("D13").Value = 1935.285702291976
$ws.Range("E13").Value = 2034.465361463813
